$d = $word.ActiveDocument

# 1. Insert the new clause about PowerShell's command color before the
#    closing parenthesis of "(which is symbolic for a full moon)".
$found = $d.Content.Find.Execute(
    "(which is symbolic for a full moon) ",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "(which is symbolic for a full moon and also matches with PowerShell’s command color) ",
    2)
Write-Output "replace1=$found"

# 2. Merge the old split around "pierce" / "s through" back together so the
#    sentence reads "...font pierces through..." as a single run (the
#    _GoBack bookmark that used to sit between them is relocated below).
$found2 = $d.Content.Find.Execute(
    "font pierces through",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "font pierces through",
    2)
Write-Output "normalize2=$found2"

# 3. Relocate the (singleton, hidden) "_GoBack" bookmark to sit right after
#    "(which is symbolic " — i.e. right before "for a full moon" — matching
#    where the author's edit actually landed.
$text = $d.Content.Text
$pos = $text.IndexOf("(which is symbolic ") + ("(which is symbolic ").Length
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)
Write-Output "bookmark moved to $pos"
